$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so that paragraph indices for
# not-yet-processed paragraphs remain stable while we insert new paragraphs.
# ---------------------------------------------------------------------------

# --- Case 3 section: fill the empty paragraph right after the "Case 3" heading
$d.Paragraphs(21).Range.Text = "The client doesn" + [char]0x2019 + "t know the end point of the server where metadata and data are stored; it queries then a well-known FHIR server to discover the endpoint fulfilling specific search criteria. This information is then used to search and access data and metadata as described for case 1 and 2."

# --- Case 2 section: fill the empty paragraph right after the "Case 2" heading,
#     then add three new paragraphs after it (two sentences + one blank line).
$d.Paragraphs(13).Range.Text = "A set of FHIR Servers are used as Data Repository by data sources. A well-known community FHIR server is used to publish appropriate FHIR resources representing metadata. "

$d.Paragraphs(13).Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.Text = "All potential clients belonging to this community know the endpoint of this well-known FHIR server."

$d.Paragraphs(14).Range.InsertParagraphAfter()
$d.Paragraphs(15).Range.Text = "Client uses FHIR API to search and get FHIR resources representing metadata; data references provided by the retrieved resources are then used to get data from the proper data repository (not known a -priori)."

$d.Paragraphs(15).Range.InsertParagraphAfter()

# --- Case 1 section: fill the empty paragraph right after the "Case 1" heading,
#     then add one new paragraph after it.
$d.Paragraphs(5).Range.Text = "All data sources store FHIR resources representing data and metadata to be shared in a well-known FHIR server. All potential clients belonging to this community know the endpoint of this well-known FHIR server."

$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "Client uses FHIR API to search and get FHIR resources representing data and/or metadata."

# --- Overview section: fill the empty paragraph after the intro sentence,
#     then add a new paragraph between the intro sentence and it.
$d.Paragraphs(3).Range.Text = "It is assumed that a proper privacy and security layer is realized to assure that all the access and usage conditions specified by the data sources are properly enforced (not shown in the pictures)."

$d.Paragraphs(2).Range.InsertParagraphAfter()
$d.Paragraphs(3).Range.Text = "This is not the only possible approach, for example hybrid FHIR/non-FHIR solutions could be considered as well."
$d.Paragraphs(3).Style = "No Spacing"

# --- Rewrite the introductory sentence itself.
$old = "This page summarizes a list of possible deployment architectures starting from a simple case where a well-known FHIR server is acting as registry and repository for data and metadata, up to the case where endpoints are supposed to be discovered. "
$new = "This page provides a short overview of a non-exhaustive list of possible logical deployment architectures that communities could consider allowing the discovery and access to FAIR data by using HL7 FHIR API; starting from a simple case where a well-known FHIR server is acting as registry and repository for data and metadata, up to the case where repository endpoints are not known a-priori and needs to be discovered."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
